$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 607
$ws.Range("H2").Value = 1800

# Row 3
$ws.Range("D3").Value = 0.135

# Row 4
$ws.Range("D4").Value = 0.105
$ws.Range("E4").Value = -0.01
$ws.Range("G4").Value = 0.093
$ws.Range("I4").Value = "0.08***"

# Row 11
$ws.Range("C11").Value = 2542
$ws.Range("D11").Value = 469
$ws.Range("F11").Value = 557
$ws.Range("G11").Value = 764
$ws.Range("H11").Value = 752

# Row 12
$ws.Range("D12").Value = 0.08
$ws.Range("F12").Value = 0.061
$ws.Range("G12").Value = 0.055
$ws.Range("H12").Value = 0.053

# Row 13
$ws.Range("C13").Value = 0.02
$ws.Range("D13").Value = 0.034
$ws.Range("F13").Value = -0.021
$ws.Range("H13").Value = 0.058
